# Updates cryptos list figures (price/volume columns) to match the
# latest scrape, and fixes row ordering for a couple of coins whose
# rank swapped (PEPE/Bittensor, BabyDogeCoin -> Optimism/Cronos/Mantle).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range('D2').NumberFormat = "@"
$ws.Range('D2').Value = '68.312.76'
$ws.Range('E2').Value = '  +1.78%  '

# Row 3
$ws.Range('D3').NumberFormat = "@"
$ws.Range('D3').Value = '2.510.50'
$ws.Range('E3').Value = '  +1.72%  '

# Row 4
$ws.Range('D4').NumberFormat = "@"
$ws.Range('D4').Value = '0.999'

# Row 5
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '592.01'
$ws.Range('E5').Value = '  +1.57%  '

# Row 6
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '176.83'
$ws.Range('E6').Value = '  +1.58%  '

# Row 7
$ws.Range('E7').Value = '  -0.05%  '

# Row 8
$ws.Range('E8').Value = '  +1.05%  '

# Row 9
$ws.Range('D9').NumberFormat = "@"
$ws.Range('D9').Value = '2.508.70'
$ws.Range('E9').Value = '  +1.66%  '

# Row 10
$ws.Range('E10').Value = '  +4.26%  '

# Row 11
$ws.Range('E11').Value = '  -1.11%  '

# Row 12
$ws.Range('E12').Value = '  +0.89%  '

# Row 13
$ws.Range('E13').Value = '  +1.11%  '

# Row 14
$ws.Range('D14').NumberFormat = "@"
$ws.Range('D14').Value = '2.997.23'
$ws.Range('E14').Value = '  +2.34%  '

# Row 15
$ws.Range('D15').NumberFormat = "@"
$ws.Range('D15').Value = '25.83'
$ws.Range('E15').Value = '  +1.88%  '

# Row 16
$ws.Range('D16').NumberFormat = "@"
$ws.Range('D16').Value = '68.057.55'
$ws.Range('E16').Value = '  +1.54%  '

# Row 17
$ws.Range('E17').Value = '  +0.39%  '

# Row 18
$ws.Range('D18').NumberFormat = "@"
$ws.Range('D18').Value = '2.495.12'
$ws.Range('E18').Value = '  +0.68%  '

# Row 19
$ws.Range('D19').NumberFormat = "@"
$ws.Range('D19').Value = '11.01'
$ws.Range('E19').Value = '  +1.07%  '

# Row 20
$ws.Range('D20').NumberFormat = "@"
$ws.Range('D20').Value = '7.44'
$ws.Range('E20').Value = '  -0.40%  '

# Row 21
$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value = '351.19'
$ws.Range('E21').Value = '  +0.74%  '

# Row 22
$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value = '4.19'
$ws.Range('E22').Value = '  +5.18%  '

# Row 23
$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value = '71.28'
$ws.Range('E23').Value = '  +2.91%  '

# Row 24
$ws.Range('E24').Value = '  -0.01%  '

# Row 25
$ws.Range('E25').Value = '  +0.37%  '

# Row 26
$ws.Range('E26').Value = '  -4.04%  '

# Row 27
$ws.Range('D27').NumberFormat = "@"
$ws.Range('D27').Value = '9.15'
$ws.Range('E27').Value = '  +0.30%  '

# Row 28
$ws.Range('E28').Value = '  +1.84%  '

# Row 29
$ws.Range('D29').NumberFormat = "@"
$ws.Range('D29').Value = '0.999'
$ws.Range('E29').Value = '  -0.44%  '

# Row 30
$ws.Range('B30').Value = 'Bittensor'
$ws.Range('C30').Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range('D30').NumberFormat = "@"
$ws.Range('D30').Value = '511.01'
$ws.Range('E30').Value = '  +2.27%  '

# Row 31
$ws.Range('B31').Value = 'PEPE'
$ws.Range('C31').Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range('D31').NumberFormat = "@"
$ws.Range('D31').Value = '0.0₃0895'
$ws.Range('E31').Value = '  -0.41%  '

# Row 32
$ws.Range('D32').NumberFormat = "@"
$ws.Range('D32').Value = '7.78'
$ws.Range('E32').Value = '  +0.80%  '

# Row 33
$ws.Range('D33').NumberFormat = "@"
$ws.Range('D33').Value = '1.26'
$ws.Range('E33').Value = '  +2.19%  '

# Row 34
$ws.Range('E34').Value = '  +1.15%  '

# Row 35
$ws.Range('E35').Value = '  -0.01%  '

# Row 37
$ws.Range('D37').NumberFormat = "@"
$ws.Range('D37').Value = '162.14'
$ws.Range('E37').Value = '  +0.37%  '

# Row 38
$ws.Range('D38').NumberFormat = "@"
$ws.Range('D38').Value = '18.69'
$ws.Range('E38').Value = '  +0.08%  '

# Row 39
$ws.Range('D39').NumberFormat = "@"
$ws.Range('D39').Value = '18.34'
$ws.Range('E39').Value = '  +1.27%  '

# Row 40
$ws.Range('E40').Value = '  -0.08%  '

# Row 41
$ws.Range('E41').Value = '  +0.01%  '

# Row 42
$ws.Range('D42').NumberFormat = "@"
$ws.Range('D42').Value = '1.75'
$ws.Range('E42').Value = '  +3.65%  '

# Row 43
$ws.Range('E43').Value = '  +0.28%  '

# Row 44
$ws.Range('E44').Value = '  +0.00%  '

# Row 45
$ws.Range('E45').Value = '  +1.89%  '

# Row 46
$ws.Range('D46').NumberFormat = "@"
$ws.Range('D46').Value = '150.88'
$ws.Range('E46').Value = '  +6.27%  '

# Row 47
$ws.Range('E47').Value = '  +2.69%  '

# Row 48
$ws.Range('E48').Value = '  +1.67%  '

# Row 49
$ws.Range('B49').Value = 'Optimism'
$ws.Range('C49').Value = 'https://coinranking.com/coin/n1p-s_gm1+optimism-op'
$ws.Range('D49').NumberFormat = "@"
$ws.Range('D49').Value = '1.60'
$ws.Range('E49').Value = '  +1.89%  '

# Row 50
$ws.Range('B50').Value = 'Cronos'
$ws.Range('C50').Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range('D50').NumberFormat = "@"
$ws.Range('D50').Value = '0.0739'
$ws.Range('E50').Value = '  +0.16%  '

# Row 51
$ws.Range('B51').Value = 'Mantle'
$ws.Range('C51').Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range('D51').NumberFormat = "@"
$ws.Range('D51').Value = '0.577'
$ws.Range('E51').Value = '  -0.70%  '

Write-Output "Updated $($ws.Name) with latest crypto figures"
